$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.6607870052619538
$ws.Cells.Item(2, 3).Value = 0.7035922067296614
$ws.Cells.Item(2, 4).Value = 0.6607870052619538
$ws.Cells.Item(2, 5).Value = 0.6594390440119021
$ws.Cells.Item(2, 6).Value = 0.6758636467627545
$ws.Cells.Item(2, 7).Value = 0.7140070347185159
$ws.Cells.Item(2, 8).Value = 0.6758636467627545
$ws.Cells.Item(2, 9).Value = 0.6751043043023429
$ws.Cells.Item(2, 10).Value = 0.6951727293525508
$ws.Cells.Item(2, 11).Value = 0.7216533663203817
$ws.Cells.Item(2, 12).Value = 0.6951727293525508
$ws.Cells.Item(2, 13).Value = 0.6942983243939361
$ws.Cells.Item(2, 14).Value = 0.7596202242049875
$ws.Cells.Item(2, 15).Value = 0.7723836279450254
$ws.Cells.Item(2, 16).Value = 0.7596202242049875
$ws.Cells.Item(2, 17).Value = 0.758967210250951
$ws.Cells.Item(2, 18).Value = 0.7617250057195151
$ws.Cells.Item(2, 19).Value = 0.7721246267043769
$ws.Cells.Item(2, 20).Value = 0.7617250057195151
$ws.Cells.Item(2, 21).Value = 0.7612488421848909
$ws.Cells.Item(2, 22).Value = 0.7768474033401969
$ws.Cells.Item(2, 23).Value = 0.7859811711427482
$ws.Cells.Item(2, 24).Value = 0.7768474033401969
$ws.Cells.Item(2, 25).Value = 0.7771016277200313
$ws.Cells.Item(3, 2).Value = 0.8068862960420956
$ws.Cells.Item(3, 3).Value = 0.8165817752422347
$ws.Cells.Item(3, 4).Value = 0.8068862960420956
$ws.Cells.Item(3, 5).Value = 0.8068985724823365
$ws.Cells.Item(3, 6).Value = 0.8111873713109128
$ws.Cells.Item(3, 7).Value = 0.8197327317326465
$ws.Cells.Item(3, 8).Value = 0.8111873713109128
$ws.Cells.Item(3, 9).Value = 0.8127405508901064
$ws.Cells.Item(3, 10).Value = 0.8240677190574239
$ws.Cells.Item(3, 11).Value = 0.8341145190108777
$ws.Cells.Item(3, 12).Value = 0.8240677190574239
$ws.Cells.Item(3, 13).Value = 0.8242075632465828
$ws.Cells.Item(3, 14).Value = 0.8477236330359185
$ws.Cells.Item(3, 15).Value = 0.8563353892598808
$ws.Cells.Item(3, 16).Value = 0.8477236330359185
$ws.Cells.Item(3, 17).Value = 0.8477534818572938
$ws.Cells.Item(3, 18).Value = 0.8498512926103867
$ws.Cells.Item(3, 19).Value = 0.8581624833208151
$ws.Cells.Item(3, 20).Value = 0.8498512926103867
$ws.Cells.Item(3, 21).Value = 0.8498447448756513
$ws.Cells.Item(3, 22).Value = 0.8498055364905056
$ws.Cells.Item(3, 23).Value = 0.8559220440159205
$ws.Cells.Item(3, 24).Value = 0.8498055364905056
$ws.Cells.Item(3, 25).Value = 0.8494250209424117
$ws.Cells.Item(4, 2).Value = 0.8433310455273393
$ws.Cells.Item(4, 3).Value = 0.8487473843105618
$ws.Cells.Item(4, 4).Value = 0.8433310455273393
$ws.Cells.Item(4, 5).Value = 0.8434160106629898
$ws.Cells.Item(4, 6).Value = 0.8712880347746511
$ws.Cells.Item(4, 7).Value = 0.8775425274847303
$ws.Cells.Item(4, 8).Value = 0.8712880347746511
$ws.Cells.Item(4, 9).Value = 0.871549669972255
$ws.Cells.Item(4, 10).Value = 0.8498055364905056
$ws.Cells.Item(4, 11).Value = 0.8564366739764255
$ws.Cells.Item(4, 12).Value = 0.8498055364905056
$ws.Cells.Item(4, 13).Value = 0.8502511645780464
$ws.Cells.Item(4, 14).Value = 0.8541523678792039
$ws.Cells.Item(4, 15).Value = 0.8589862301481432
$ws.Cells.Item(4, 16).Value = 0.8541523678792039
$ws.Cells.Item(4, 17).Value = 0.853473523710947
$ws.Cells.Item(4, 18).Value = 0.8584534431480211
$ws.Cells.Item(4, 19).Value = 0.8632781427252849
$ws.Cells.Item(4, 20).Value = 0.8584534431480211
$ws.Cells.Item(4, 21).Value = 0.8578926526236801
$ws.Cells.Item(4, 22).Value = 0.8563029055136125
$ws.Cells.Item(4, 23).Value = 0.8610149354238773
$ws.Cells.Item(4, 24).Value = 0.8563029055136125
$ws.Cells.Item(4, 25).Value = 0.8556654499841484
$ws.Cells.Item(5, 2).Value = 0.8348204072294669
$ws.Cells.Item(5, 3).Value = 0.8424913763141845
$ws.Cells.Item(5, 4).Value = 0.8348204072294669
$ws.Cells.Item(5, 5).Value = 0.8328226448728516
$ws.Cells.Item(5, 6).Value = 0.8498055364905056
$ws.Cells.Item(5, 7).Value = 0.8569006817846342
$ws.Cells.Item(5, 8).Value = 0.8498055364905056
$ws.Cells.Item(5, 9).Value = 0.8486362112405953
$ws.Cells.Item(5, 14).Value = 0.8347746511095859
$ws.Cells.Item(5, 15).Value = 0.8453118490272775
$ws.Cells.Item(5, 16).Value = 0.8347746511095859
$ws.Cells.Item(5, 17).Value = 0.8329856798001873
$ws.Cells.Item(5, 18).Value = 0.8412262640128118
$ws.Cells.Item(5, 19).Value = 0.8519692035217915
$ws.Cells.Item(5, 20).Value = 0.8412262640128118
$ws.Cells.Item(5, 21).Value = 0.8394283439992805
$ws.Cells.Item(6, 2).Value = 0.8541066117593228
$ws.Cells.Item(6, 3).Value = 0.859921301250685
$ws.Cells.Item(6, 4).Value = 0.8541066117593228
$ws.Cells.Item(6, 5).Value = 0.8540234574043056
$ws.Cells.Item(6, 6).Value = 0.8755891100434683
$ws.Cells.Item(6, 7).Value = 0.8809666747163373
$ws.Cells.Item(6, 8).Value = 0.8755891100434683
$ws.Cells.Item(6, 9).Value = 0.8752460492178475
$ws.Cells.Item(6, 10).Value = 0.8219171814230153
$ws.Cells.Item(6, 11).Value = 0.8347579140389392
$ws.Cells.Item(6, 12).Value = 0.8219171814230153
$ws.Cells.Item(6, 13).Value = 0.8227163881853119
$ws.Cells.Item(6, 14).Value = 0.8540837336993823
$ws.Cells.Item(6, 15).Value = 0.8592373894203723
$ws.Cells.Item(6, 16).Value = 0.8540837336993823
$ws.Cells.Item(6, 17).Value = 0.85340309945919
$ws.Cells.Item(6, 18).Value = 0.8519103180050331
$ws.Cells.Item(6, 19).Value = 0.8594872292154676
$ws.Cells.Item(6, 20).Value = 0.8519103180050331
$ws.Cells.Item(6, 21).Value = 0.8510956953545248
$ws.Cells.Item(6, 22).Value = 0.8519103180050331
$ws.Cells.Item(6, 23).Value = 0.8563462984869019
$ws.Cells.Item(6, 24).Value = 0.8519103180050331
$ws.Cells.Item(6, 25).Value = 0.851092372017099